# Calibration data re-sort: the data rows (A2:D8) need to end up sorted
# in ascending order of column A (time). Row 2 already holds the smallest
# time value; rows 3-8 get reordered around it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRows = 7
$numCols = 4

# Read the current data block (rows 2-8, columns A-D) into memory.
$data = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $ws.Cells.Item($r + 2, $c + 1).Value2
    }
}

# Pair each row's time (column A) with its original row index, then sort
# those pairs by time ascending.
$pairs = @()
for ($r = 0; $r -lt $numRows; $r++) {
    $pairs += @{ Key = $data[$r, 0]; Row = $r }
}
$sortedPairs = $pairs | Sort-Object { $_.Key }

# Write the rows back out in the new (sorted) order.
for ($i = 0; $i -lt $numRows; $i++) {
    $srcRow = $sortedPairs[$i].Row
    $destRow = $i + 2
    for ($c = 0; $c -lt $numCols; $c++) {
        $ws.Cells.Item($destRow, $c + 1).Value2 = $data[$srcRow, $c]
    }
}
